# Weekly data refresh: a new sample week is inserted at the top of the
# recent data block (row 12), pushing the existing rows 12-23 down to
# rows 13-24. The sheet grows from A1:R23 to A1:R24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 12; this shifts rows
# 12-23 down to 13-24 (and extends the used range to A1:R24).
$ws.Rows(12).Insert()

# Populate the newly inserted row 12 with this week's record.
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 45079
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 100112013
$ws.Cells.Item(12, 7).Value = "Alcachofa"
$ws.Cells.Item(12, 8).Value = "Madrigal"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 11).Value = 16500
$ws.Cells.Item(12, 12).Value = 17000
$ws.Cells.Item(12, 13).Value = 16750
$ws.Cells.Item(12, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(12, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 16).Value = 419
$ws.Cells.Item(12, 17).Value = 40
$ws.Cells.Item(12, 18).Value = "Hortaliza"
